# Updated cryptos list on Sun Aug 18 04:39:31 UTC 2024 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) cells are stored as
# plain text in the workbook (t="inlineStr"), even when the text looks like
# a plain number (e.g. "535.53"). Assigning such a string straight to
# Range.Value lets the COM layer auto-coerce it to a numeric cell, which
# would change the cell's stored type. To keep these cells as text - exactly
# like the source file - each Price cell is (1) explicitly formatted as
# Text, (2) written, then (3) reset to the "Normal" style so no stray
# number-format / style index is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "59.356.05"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.598.32"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "535.53"
$ws.Range("E5").Value = "  +2.46%  "

# Row 6 - Solana
Set-TextValue "D6" "141.11"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.15%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.32%  "

# Row 9 - Toncoin
Set-TextValue "D9" "6.48"
$ws.Range("E9").Value = "  -1.34%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.20%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.334"
$ws.Range("E11").Value = "  +1.55%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.83%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.061.49"
$ws.Range("E13").Value = "  +0.65%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "59.276.01"
$ws.Range("E14").Value = "  +0.54%  "

# Row 15 - Avalanche
Set-TextValue "D15" "20.70"
$ws.Range("E15").Value = "  +1.21%  "

# Row 16 - swapped to WrappedEther (was ShibaInu)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.604.84"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17 - swapped to ShibaInu (was WrappedEther)
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000133"
$ws.Range("E17").Value = "  +0.09%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "340.85"
$ws.Range("E18").Value = "  +0.62%  "

# Row 19 - Polkadot
Set-TextValue "D19" "4.36"
$ws.Range("E19").Value = "  +1.43%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.17%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.46%  "

# Row 22 - Dai
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Litecoin
Set-TextValue "D23" "67.51"
$ws.Range("E23").Value = "  +2.30%  "

# Row 24 - Polygon
$ws.Range("E24").Value = "  +1.29%  "

# Row 25 - Kaspa
Set-TextValue "D25" "0.165"
$ws.Range("E25").Value = "  -1.51%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.18%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "7.22"
$ws.Range("E27").Value = "  +2.97%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +2.23%  "

# Row 29 - USDe
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +4.70%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -1.84%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "18.81"
$ws.Range("E32").Value = "  +0.62%  "

# Row 33 - Monero
Set-TextValue "D33" "149.81"
$ws.Range("E33").Value = "  +0.56%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  -0.54%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -1.22%  "

# Row 36 - Stacks
$ws.Range("E36").Value = "  -0.75%  "

# Row 37 - SuiNetwork
$ws.Range("E37").Value = "  +1.52%  "

# Row 38 - Fetch.AI
Set-TextValue "D38" "0.822"
$ws.Range("E38").Value = "  -0.68%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  +0.41%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.24%  "

# Row 41 - Bittensor
Set-TextValue "D41" "272.74"
$ws.Range("E41").Value = "  +0.62%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.597"
$ws.Range("E42").Value = "  +1.55%  "

# Row 43 - WhiteBITCoin
Set-TextValue "D43" "10.74"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44 - Stellar
Set-TextValue "D44" "0.0952"
$ws.Range("E44").Value = "  -0.19%  "

# Row 45 - Hedera
$ws.Range("E45").Value = "  +1.12%  "

# Row 46 - InjectiveProtocol
Set-TextValue "D46" "18.57"
$ws.Range("E46").Value = "  +3.36%  "

# Row 47 - Maker
Set-TextValue "D47" "1.941.21"
$ws.Range("E47").Value = "  -1.31%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +1.34%  "

# Row 49 - RenderToken
Set-TextValue "D49" "4.49"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50 - Aave
Set-TextValue "D50" "111.58"
$ws.Range("E50").Value = "  -1.55%  "

# Row 51 - ZEEBU
Set-TextValue "D51" "4.74"
$ws.Range("E51").Value = "  +0.21%  "
